# #4 ajout du numero de la release/milestone pour chaque slice
# Fill column F (release/milestone number) for rows 17-47 of the "Feuil1"
# tracking sheet, then leave the selection where the author left it (F47).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

$values = @{
    17 = 4
    18 = 4
    19 = 4
    20 = 4
    21 = 4
    22 = 4
    23 = 5
    24 = 5
    25 = 5
    26 = 5
    27 = 6
    28 = 6
    29 = 6
    30 = 6
    31 = 6
    32 = 6
    33 = 6
    34 = 6
    35 = 7
    36 = 7
    37 = 7
    38 = 7
    39 = 7
    40 = 7
    41 = 7
    42 = 8
    43 = 8
    44 = 8
    45 = 8
    46 = 8
    47 = 9
}

foreach ($row in $values.Keys) {
    # Column F is column 6
    $ws.Cells.Item($row, 6).Value = $values[$row]
}

# Match the saved view state: scrolled so row 47 is visible, with F47 selected.
$ws.Range("F47").Select()
